$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$pQuestions = '<w:p><w:r><w:rPr><w:i/></w:rPr><w:t>What is</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>/are</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> your </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>research question(s)</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>?</w:t></w:r><w:r><w:t xml:space="preserve"> (Answer this after you complete the following worksheet).</w:t></w:r></w:p>'
$pStep1 = '<w:p><w:r><w:t xml:space="preserve">Step 1: Visit the Wikipedia page for your topic. You may have to modify your topic or controversy to find a page that provides actual information related to the subject that interests you. </w:t></w:r></w:p>'
$pMla = '<w:p><w:r><w:t>Write the MLA citation for the page. You can find out how to write an MLA citation for Wikipedia by googling it. (Hint: it has its own Wikipedia page).</w:t></w:r></w:p>'
$pStep2 = '<w:p><w:r><w:t>Step 2: Read the Wikipedia page</w:t></w:r><w:r><w:t>, focusing on relevant sections</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> Are there any explicitly stated controversies (e.g.: “People disagree about x…” or “Some people think x, while others think y…”)? List them here:</w:t></w:r></w:p>'
$pStep3 = '<w:p><w:r><w:t>Step 3: Does the Wikipedia article mention any relevant current events (e.g.: a court case; a controversial article, film, or book; new legislation; a presidential speech; etc.)? List them here:</w:t></w:r></w:p>'
$pStep4 = '<w:p><w:r><w:t>Step 4</w:t></w:r><w:r><w:t xml:space="preserve">: Based on what you have read, </w:t></w:r><w:r><w:t>list at least four search terms you want to follow up with. They might be names of people or places related to your controversy; specific events that happened related to your controversy; jargon that is used in talking about your topic; etc.</w:t></w:r></w:p>'
$pStep5 = '<w:p><w:r><w:t>Step 5</w:t></w:r><w:r><w:t>: Look at the citations on the bottom of the page. Copy and paste any useful-looking links here:</w:t></w:r></w:p>'
$pStep6 = '<w:p><w:r><w:t>Step 6: Fill in your research questions at the top.</w:t></w:r></w:p>'
$pStep7 = '<w:p><w:r><w:t xml:space="preserve">Step 7: Where are you now? Do you have too much information, or not enough? If you have too much information, how are you going to narrow your topic so that you can filter? Consider specific subcategories that catch your attention. If you don’t have enough information, how might you modify your topic to be more effective? Consider new keywords, new Wikipedia pages, or other sources that might serve you better.  </w:t></w:r></w:p>'
$pRepeat = '<w:p><w:r><w:t>Repeat this process as necessary if you find other helpful Wikipedia pages.</w:t></w:r></w:p>'
$pStay = '<w:p><w:r><w:t>Stay focused: don’t go down the Wikipedia rabbit hole! You can lose days down there.</w:t></w:r></w:p>'

# ---------------------------------------------------------------------------
# rhe306-spring2014/documents/wikipedia.docx edit:
#  - Reword the "What is your controversy?" prompt into a multi-run
#    "What is/are your research question(s)?" prompt followed by new
#    instruction text.
#  - Insert a brand-new worksheet (Steps 1-7, "Repeat this process...",
#    "Stay focused...") directly after that prompt, built from the old
#    Step 1-5 paragraphs (unchanged) plus new Step 6/Step 7 copy and two
#    new closing paragraphs.
#  - Remove the old duplicate Step 1-6 block that used to follow the
#    _GoBack bookmark paragraph at the end of the document.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# Paragraph 5 is "What is your controversy? (If you don't know, come back
# to this question at the end of this assignment)." -- replace its whole
# range (this also drops the old paragraph mark's rsid/paraId, matching a
# freshly authored <w:p>).
$targetPara = $d.Paragraphs.Item(5)
$targetRange = $targetPara.Range

# Build the replacement block: the reworded question paragraph, then the
# full new worksheet, each entry followed by a blank paragraph, exactly as
# laid out in the target document.
$block = $pQuestions + '<w:p/>'
$block += $pStep1 + '<w:p/>'
$block += $pMla + '<w:p/>' + '<w:p/>'
$block += $pStep2 + '<w:p/>'
$block += $pStep3 + '<w:p/>'
$block += $pStep4 + '<w:p/>'
$block += $pStep5 + '<w:p/>'
$block += $pStep6 + '<w:p/>'
$block += $pStep7 + '<w:p/>'
$block += $pRepeat + '<w:p/>'
$block += $pStay + '<w:p/>'

# Namespace only needs to be declared once, on the first element of the
# fragment; it is inherited by the remaining sibling elements.
$block = $block -replace '<w:p>', ('<w:p ' + $wNs + '>'), 1

$targetRange.InsertXML($block)

# The old duplicate "Step 1" ... "Step 6" block (unchanged content) still
# sits right after the _GoBack bookmark paragraph, further down in the
# document; the new block above already reproduces Steps 1-5 verbatim and
# replaces Step 6 with fresh copy, so delete the leftover duplicate run of
# paragraphs up to (and including) the final empty paragraph before
# </w:body>.
#
# Locate the duplicate block precisely: it starts with the paragraph whose
# text begins "Step 1: Visit the Wikipedia page" and which comes *after*
# the newly-inserted worksheet (i.e. the second occurrence in the
# document), and it runs through to the last paragraph of the body.
$stepOneIndices = @()
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith('Step 1: Visit the Wikipedia page')) {
        $stepOneIndices += $i
    }
}

if ($stepOneIndices.Count -ge 2) {
    $dupStart = $stepOneIndices[1]
    $lastParaIndex = $d.Paragraphs.Count
    $startRange = $d.Paragraphs.Item($dupStart).Range
    $endRange = $d.Paragraphs.Item($lastParaIndex).Range
    $fullDelete = $d.Range($startRange.Start, $endRange.End)
    $fullDelete.Delete()
}
